$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()
$ws.Range("J5").Activate()
